$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText($ref, $val) {
    $r = $ws.Range($ref)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.NumberFormat = "General"
    $r.Style = "Normal"
}

Set-CellText "D2" "26.656.79"
Set-CellText "E2" "  -0.17%  "
Set-CellText "D3" "1.529.65"
Set-CellText "E3" "  -1.88%  "
Set-CellText "E4" "  -0.15%  "
Set-CellText "D5" "205.31"
Set-CellText "E5" "  -0.16%  "
Set-CellText "E6" "  -0.90%  "
Set-CellText "E7" "  -0.04%  "
Set-CellText "E8" "  -2.87%  "
Set-CellText "E9" "  -1.28%  "
Set-CellText "E10" "  -0.80%  "
Set-CellText "E11" "  -1.30%  "
Set-CellText "D12" "1.747.15"
Set-CellText "E12" "  -1.92%  "
Set-CellText "D13" "1.532.37"
Set-CellText "E13" "  -1.65%  "
Set-CellText "E14" "  -2.11%  "
Set-CellText "E15" "  -1.01%  "
Set-CellText "D16" "61.39"
Set-CellText "D17" "26.656.96"
Set-CellText "E17" "  -0.30%  "
Set-CellText "D18" "212.35"
Set-CellText "E18" "  -0.76%  "
Set-CellText "E19" "  +1.22%  "
Set-CellText "E20" "  -2.16%  "
Set-CellText "E21" "  -0.06%  "
Set-CellText "E22" "  -2.26%  "
Set-CellText "E23" "  -3.39%  "
Set-CellText "E24" "  -3.31%  "
Set-CellText "D25" "151.92"
Set-CellText "E25" "  -0.72%  "
Set-CellText "D26" "6.51"
Set-CellText "E26" "  -3.71%  "
Set-CellText "E27" "  +0.01%  "
Set-CellText "E28" "  -0.16%  "
Set-CellText "E29" "  -0.78%  "
Set-CellText "D30" "1.10"
Set-CellText "E30" "  -0.97%  "
Set-CellText "E31" "  -1.88%  "
Set-CellText "E32" "  +2.72%  "
Set-CellText "D33" "1.351.24"
Set-CellText "E33" "  -2.54%  "
Set-CellText "E34" "  -0.02%  "
Set-CellText "E35" "  -3.73%  "
Set-CellText "D36" "0.951"
Set-CellText "E36" "  +2.27%  "
Set-CellText "D37" "2.27"
Set-CellText "E37" "  -0.65%  "
Set-CellText "E38" "  +0.17%  "
Set-CellText "E39" "  +0.77%  "
Set-CellText "B40" "ARBITRUM"
Set-CellText "C40" "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-CellText "D40" "0.796"
Set-CellText "E40" "  -1.46%  "
Set-CellText "B41" "FraxShare"
Set-CellText "C41" "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-CellText "D41" "5.69"
Set-CellText "E41" "  +5.19%  "
Set-CellText "D42" "0.994"
Set-CellText "E42" "  +0.30%  "
Set-CellText "E43" "  -0.09%  "
Set-CellText "E44" "  -1.28%  "
Set-CellText "E45" "  -2.24%  "
Set-CellText "D46" "2.26"
Set-CellText "E46" "  -3.51%  "
Set-CellText "D47" "1.662.05"
Set-CellText "E47" "  -1.96%  "
Set-CellText "D48" "85.54"
Set-CellText "E48" "  -0.02%  "
Set-CellText "D49" "0.0506"
Set-CellText "E49" "  +2.78%  "
Set-CellText "D50" "0.0₇0963"
Set-CellText "E50" "  -2.21%  "
Set-CellText "E51" "  -0.17%  "
